# Restore C10 ("From" value of rule R20 on the Rules sheet) back to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
